$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-27 Sunday", "2024-10-28 Monday"),
    @("89÷5=", "81÷4="),
    @("37÷2=", "87÷4="),
    @("28÷8=", "92÷2="),
    @("51÷5=", "87÷7="),
    @("40÷5=", "32÷4="),
    @("41÷2=", "64÷8="),
    @("79÷4=", "58÷7="),
    @("57÷4=", "51÷7="),
    @("74÷9=", "20÷9="),
    @("31÷8=", "49÷4="),
    @("44÷7=", "42÷5="),
    @("42÷4=", "75÷7="),
    @("64÷5=", "25÷2="),
    @("86÷2=", "51÷8="),
    @("37÷8=", "82÷2="),
    @("80÷4=", "48÷4="),
    @("44÷9=", "76÷9="),
    @("60÷7=", "98÷6="),
    @("93÷4=", "91÷7="),
    @("26÷3=", "14÷3="),
    @("20÷4=", "56÷7="),
    @("26÷6=", "77÷6="),
    @("52÷8=", "16÷5="),
    @("12÷9=", "30÷9="),
    @("68÷4=", "18÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
